$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a string value to a cell while forcing Text format so
# numeric-looking strings (e.g. '1.00', '3.20') are preserved verbatim
# instead of being parsed into Excel numbers, and restore the cell's
# original style afterwards so no stray formatting is introduced.
function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

$updates = @(
    @{ Row = 2; D = '87.025.66'; E = '  +9.74%  ' },
    @{ Row = 3; D = '3.355.38'; E = '  +6.31%  ' },
    @{ Row = 4; E = '  -0.26%  ' },
    @{ Row = 5; D = '221.60'; E = '  +8.09%  ' },
    @{ Row = 6; D = '640.51'; E = '  +2.68%  ' },
    @{ Row = 7; D = '0.325'; E = '  +22.86%  ' },
    @{ Row = 8; D = '0.998'; E = '  -0.25%  ' },
    @{ Row = 9; D = '0.628'; E = '  +6.77%  ' },
    @{ Row = 10; D = '3.360.21'; E = '  +6.30%  ' },
    @{ Row = 11; D = '0.616'; E = '  +4.29%  ' },
    @{ Row = 12; D = '0.0000276'; E = '  +9.30%  ' },
    @{ Row = 13; E = '  +1.96%  ' },
    @{ Row = 14; D = '3.958.88'; E = '  +5.88%  ' },
    @{ Row = 15; D = '34.60'; E = '  +10.44%  ' },
    @{ Row = 16; D = '5.42'; E = '  +2.96%  ' },
    @{ Row = 17; D = '86.764.84'; E = '  +9.21%  ' },
    @{ Row = 18; D = '3.334.56'; E = '  +5.48%  ' },
    @{ Row = 19; D = '14.80'; E = '  +3.77%  ' },
    @{ Row = 20; D = '3.22'; E = '  +10.52%  ' },
    @{ Row = 21; D = '449.59'; E = '  +3.22%  ' },
    @{ Row = 22; D = '9.23'; E = '  +1.62%  ' },
    @{ Row = 23; D = '5.34'; E = '  +2.26%  ' },
    @{ Row = 24; E = '  +9.68%  ' },
    @{ Row = 25; D = '5.45'; E = '  +17.08%  ' },
    @{ Row = 26; D = '12.27'; E = '  +13.41%  ' },
    @{ Row = 27; D = '3.512.20'; E = '  +5.63%  ' },
    @{ Row = 28; D = '78.76'; E = '  +3.85%  ' },
    @{ Row = 29; D = '0.0000134'; E = '  +10.80%  ' },
    @{ Row = 30; D = '1.00'; E = '  +0.50%  ' },
    @{ Row = 31; D = '0.181'; E = '  +47.77%  ' },
    @{ Row = 32; D = '608.55'; E = '  +11.33%  ' },
    @{ Row = 33; D = '9.37'; E = '  +4.36%  ' },
    @{ Row = 34; E = '  +0.34%  ' },
    @{ Row = 35; D = '1.57'; E = '  +6.64%  ' },
    @{ Row = 36; E = '  +3.81%  ' },
    @{ Row = 37; D = '0.152'; E = '  +2.80%  ' },
    @{ Row = 38; D = '23.61'; E = '  +2.67%  ' },
    @{ Row = 39; D = '6.60'; E = '  +17.72%  ' },
    @{ Row = 40; D = '0.422'; E = '  +4.37%  ' },
    @{ Row = 41; B = 'dogwifhat'; C = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D = '3.20'; E = '  +20.27%  ' },
    @{ Row = 42; B = 'FirstDigitalUSD'; C = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D = '0.998'; E = '  -0.27%  ' },
    @{ Row = 43; B = 'Stacks'; C = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D = '2.09'; E = '  +16.87%  ' },
    @{ Row = 44; B = 'WhiteBITCoin'; C = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; D = '21.34'; E = '  +2.85%  ' },
    @{ Row = 45; E = '  +0.05%  ' },
    @{ Row = 46; B = 'Aave'; C = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; D = '191.34'; E = '  +2.28%  ' },
    @{ Row = 47; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '156.65'; E = '  -4.55%  ' },
    @{ Row = 48; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '1.38'; E = '  +7.02%  ' },
    @{ Row = 49; B = 'OKB'; C = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D = '45.76'; E = '  +5.79%  ' },
    @{ Row = 50; D = '0.795'; E = '  +1.99%  ' },
    @{ Row = 51; B = 'ARBITRUM'; C = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'; D = '0.665'; E = '  +6.40%  ' }
)

foreach ($u in $updates) {
    if ($u.ContainsKey('B')) { $ws.Cells.Item($u.Row, 2).Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Cells.Item($u.Row, 3).Value = $u.C }
    if ($u.ContainsKey('D')) { Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D }
    if ($u.ContainsKey('E')) { Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E }
}
